$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Tighten row 5 values down to 2 decimal places (custom accuracy rounding)
$ws.Range("B5").Value = 20.08
$ws.Range("C5").Value = 14.91
$ws.Range("D5").Value = 1.22
$ws.Range("E5").Value = 43.9
$ws.Range("F5").Value = 35.65
$ws.Range("G5").Value = 15.74
$ws.Range("H5").Value = 61.61
$ws.Range("I5").Value = 24.43
$ws.Range("J5").Value = 10.81
$ws.Range("K5").Value = 15.95
$ws.Range("L5").Value = 17.59
$ws.Range("M5").Value = 18.71
$ws.Range("N5").Value = 5.07
$ws.Range("O5").Value = 15.79
$ws.Range("P5").Value = 22.44
$ws.Range("Q5").Value = 13.37
$ws.Range("R5").Value = 0.79
$ws.Range("S5").Value = 0.83
$ws.Range("T5").Value = 233.03
$ws.Range("U5").Value = 44.16
$ws.Range("V5").Value = 14.58
$ws.Range("W5").Value = 29.62
$ws.Range("X5").Value = 15.52
$ws.Range("Y5").Value = 2.37
$ws.Range("Z5").Value = 29.96
$ws.Range("AA5").Value = 12.87
$ws.Range("AB5").Value = 11.44
$ws.Range("AC5").Value = 13.45
$ws.Range("AD5").Value = 18.45
$ws.Range("AE5").Value = 0.54
$ws.Range("AF5").Value = 56.05
$ws.Range("AG5").Value = 8.18
$ws.Range("AH5").Value = 18.22

# Remove the now-superseded last row (row 6) entirely
$ws.Rows.Item(6).Delete()
